$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6,3).Value = "[
`"Explain the value of voice`",
`"Design the user experience`",
`"Design the architecture to build the skill`",
`"Follow AWS and Alexa security best practices for the skill`",
`"Develop, test, validate, and troubleshoot the skill`",
`"Manage the skill-publishing process and work with the Alexa Developer Console`",
`"Manage skill operations and life cycles`"
]"
$ws.Rows.Item(6).RowHeight = 15
$ws.Cells.Item(1,3).Value = "descriptions"
$ws.Cells.Item(2,3).Value = "[
`"Define what the AWS Cloud is and the basic global infrastructure`",
`"Describe basic AWS Cloud architectural principles`",
`"Describe the AWS Cloud value proposition`",
`"Describe key services on the AWS platform and their common use cases (for example, compute and analytics)`",
`"Describe basic security and compliance aspects of the AWS platform and the shared security model`",
`"Define the billing, account management, and pricing models`",
`"Identify sources of documentation or technical assistance (for example, whitepapers or support tickets)`",
`"Describe basic/core characteristics of deploying and operating in the AWS Cloud`"
]"
$ws.Rows.Item(2).RowHeight = 15
$ws.Cells.Item(3,3).Value = "[
`"Demonstrate an understanding of core AWS services, uses, and basic AWS architecture best practices`",
`"Demonstrate proficiency in developing, deploying, and debugging cloud-based applications using AWS`"
]"
$ws.Rows.Item(3).RowHeight = 15
$ws.Cells.Item(4,3).Value = "[
`"Effectively demonstrate knowledge of how to architect and deploy secure and robust applications on AWS technologies`",
`"Define a solution using architectural design principles based on customer requirements`",
`"Provide implementation guidance based on best practices to the organization throughout the life cycle of the project`"
]"
$ws.Rows.Item(4).RowHeight = 15
$ws.Cells.Item(5,3).Value = "[
`"Deploy, manage, and operate scalable, highly available, and fault-tolerant systems on AWS`",
`"Implement and control the flow of data to and from AWS`",
`"Select the appropriate AWS service based on compute, data, or security requirements`",
`"Identify appropriate use of AWS operational best practices`",
`"Estimate AWS usage costs and identify operational cost control mechanisms`",
`"Migrate on-premises workloads to AWS`"
]"
$ws.Rows.Item(5).RowHeight = 15
$ws.Cells.Item(7,3).Value = "[
`"Understand and differentiate the key features of AWS database services`",
`"Analyze needs and requirements to recommend and design appropriate database solutions using AWS services`"
]"
$ws.Rows.Item(7).RowHeight = 15
$ws.Cells.Item(8,3).Value = "[
`"An understanding of specialized data classifications and AWS data protection mechanisms`",
`"An understanding of data encryption methods and AWS mechanisms to implement them`",
`"An understanding of secure Internet protocols and AWS mechanisms to implement them`",
`"A working knowledge of AWS security services and features of services to provide a secure production environment`",
`"Competency gained from two or more years of production deployment experience using AWS security services and features`",
`"Ability to make tradeoff decisions with regard to cost, security, and deployment complexity given a set of application requirements`",
`"An understanding of security operations and risk`"
]"
$ws.Rows.Item(8).RowHeight = 15
$ws.Cells.Item(9,3).Value = "[
`"Implement and manage continuous delivery systems and methodologies on AWS`",
`"Implement and automate security controls, governance processes, and compliance validation`",
`"Define and deploy monitoring, metrics, and logging systems on AWS`",
`"Implement systems that are highly available, scalable, and self-healing on the AWS platform`",
`"Design, manage, and maintain tools to automate operational processes`"
]"
$ws.Rows.Item(9).RowHeight = 15
$ws.Cells.Item(10,3).Value = "[
`"Design and deploy dynamically scalable, highly available, fault-tolerant, and reliable applications on AWS`",
`"Select appropriate AWS services to design and deploy an application based on given requirements`",
`"Migrate complex, multi-tier applications on AWS`",
`"Design and deploy enterprise-wide scalable operations on AWS`",
`"Implement cost-control strategies`"
]"
$ws.Rows.Item(10).RowHeight = 15
$ws.Cells.Item(11,3).Value = "[
`"Describe cloud concepts`",
`"Describe core Azure services`",
`"Describe core solutions and management tools on Azure`",
`"Describe general security and network security features`",
`"Describe identity, governance, privacy, and compliance features`",
`"Describe Azure cost management and Service Level Agreements`"
]"
$ws.Rows.Item(11).RowHeight = 15
$ws.Cells.Item(12,3).Value = "[
`"I/O and NIO`",
`"Generics and Collections`",
`"Flow Control and Exceptions`",
`"Strings, Arrays, and ArrayLists`",
`"Declarations and Access Control`",
`"Advanced OO and Design Patterns`",
`"Assertions and Java SE 7 Exceptions`",
`"Threads, Inner Classes, and Concurrency`",
`"String Processing, Data Formatting, and Resources Bundles`"
]"
$ws.Rows.Item(12).RowHeight = 15
$ws.Cells.Item(13,3).Value = "[
`"Explain SAFe Agile Principles`",
`"Plan Iterations`",
`"Plan Program Increments`",
`"Execute Iterations and demonstrate value`",
`"Improve Agile Release Train processes`",
`"Integrate and work with other teams on the Agile Release Train`",
`"Perform as member of an Agile Team on an Agile Release Train`"
]"
$ws.Rows.Item(13).RowHeight = 15

$ws.Range("C14").Select()
